# Updates the cryptocurrency price/volume table on Sheet1 with refreshed
# values (and, for rows 44-45, swaps the Mantle/EnergySwap entries).
# NumberFormat is forced to Text ("@") before writing the Price/Volume
# columns so that values such as "554.79" or "62.873.88" are stored as
# literal text instead of being auto-coerced into numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.873.88'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.35%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.683.99'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.88%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '554.79'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.60%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.90'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.92%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.32%  '

# Row 9
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.88%  '

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.30%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.368'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.76%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.41'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.66%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.159.67'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.84%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.55'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.38%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '62.823.49'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.19%  '

# Row 16
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.69%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.685.50'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.98%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.86'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -4.08%  '

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.62%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '345.06'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.44%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.22'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.39%  '

# Row 22
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.02%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.506'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.00%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.20'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.42%  '

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.40%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.12%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.17'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.31%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.42'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +8.01%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0856'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.40%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.24'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.39%  '

# Row 31
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.19%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '164.14'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.31%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.90'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.52%  '

# Row 34
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.19%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.48'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.84%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.78'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.32%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '348.36'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.13%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.25'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.13%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.944'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.61%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.98'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.59%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.34'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.14%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.83'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.63%  '

# Row 44
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.618'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.83%  '

# Row 45
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.17'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.30%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0558'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.07%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.01%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.00'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.50%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0970'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.01%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0241'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.39%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '128.60'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.19%  '
